$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8011
$wsExhibit.Range("F17").Value = 5907
$wsExhibit.Range("F19").Value = 275
$wsExhibit.Range("F20").Value = 1910
$wsExhibit.Range("F21").Value = 16

# Sheet "全部类型" (sheet4 / rId4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8011
$wsAll.Range("F18").Value = 5907
$wsAll.Range("F21").Value = 275
$wsAll.Range("F22").Value = 1910
$wsAll.Range("F23").Value = 16
